$d = $word.ActiveDocument

# Update the process reference number (caratula)
$d.Content.Find.Execute("DEAJGCC23-13186", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DEAJGCC23-13204", 2)

# Update all occurrences of the date
$d.Content.Find.Execute("28 de agosto de 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "02 de septiembre de 2024", 2)
